$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force Text format on the data range so numeric-looking strings
# (e.g. "240.39", "1.0000") are preserved exactly as text, matching
# the original inline-string cells, rather than being parsed as numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.400.41'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.847.48'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '240.39'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = '0.6328'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('B8').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C8').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D8').Value = '6.573.25'
$ws.Range('E8').Value = '  +207.31%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '47.22'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.07566'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').Value = '0.2963'
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '24.52'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07713'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.833.11'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '4.998'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.6858'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001010'
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = '83.12'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '6.164'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '29.423.22'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '230.19'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '12.46'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.9995'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '7.572'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = '1.000'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '156.75'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1401'
$ws.Range('E27').Value = '  +0.75%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '8.383'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '17.67'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.466'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.05743'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.250'
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.132'
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '4.028'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.849'
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.157'
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.7168'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.588'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.251.75'
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01820'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.781'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.9076'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '6.197'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '2.001.93'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '101.69'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '66.22'
$ws.Range('E47').Value = '  -1.86%  '
$ws.Range('D48').Value = '9.171'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '7.051'
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '0.4033'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.707'
$ws.Range('E51').Value = '  +1.16%  '

# Restore default (Normal) style so we do not leave stray per-cell
# number-format overrides behind.
$dataRange.Style = "Normal"

